$d = $word.ActiveDocument

# --- 1. Split the title "Response Variable Methods " into two runs and
#        retitle the first part to "Snorkel Survey", leaving " Methods " intact.
$titleRng = $d.Range(0, 17)            # "Response Variable" (17 chars)
$titleRng.Bold = 1                      # force a run split (different formatting)
$titleRng.Text = "Snorkel Survey"       # replace the text in that run

$titleRng2 = $d.Range(0, 14)            # "Snorkel Survey" (14 chars)
$titleRng2.Bold = 0                     # restore formatting so rPr matches sibling run

# --- 2. Split the first body paragraph's run in two (no text changes) right
#        before "A - 1)" so "...Appendix A (Table " / "A - 1), which..." land
#        in separate runs, matching the target markup.
$findRng = $d.Content
$findRng.Find.Execute("A - 1), which can be viewed", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $findRng.Start
$bodyStart = $d.Paragraphs(2).Range.Start

$bodyRng = $d.Range($bodyStart, $splitPoint)
$bodyRng.Bold = 1
$bodyRng.Bold = 0

# --- 3. Delete everything from the start of paragraph 3 ("Explanatory
#        Variable Methods" heading) through the end of the document, which
#        removes that whole section (heading, blank paragraph, explanatory
#        paragraph + bookmark, trailing blank paragraph).
$delStart = $d.Paragraphs(3).Range.Start
$delEnd = $d.Paragraphs($d.Paragraphs.Count).Range.End
$delRng = $d.Range($delStart, $delEnd)
$delRng.Delete()
